# Auto-generated Excel COM-interop script
# Updates numeric values across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled pricing-data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 486.875
$ws.Range("J13").Value = 156.42857
$ws.Range("L13").Value = 156.42857
$ws.Range("N13").Value = -494.42857
$ws.Range("H17").Value = 57954.25
$ws.Range("J17").Value = 57954.25
$ws.Range("L17").Value = 173862.75
$ws.Range("N17").Value = -174198.75
$ws.Range("H97").Value = 2500
$ws.Range("J97").Value = 2500
$ws.Range("L97").Value = 7500
$ws.Range("N97").Value = -8492
$ws.Range("H98").Value = 1054.6154
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 3511.5293
$ws.Range("I100").Value = 2954.4546
$ws.Range("K100").Value = 2954.4546
$ws.Range("M100").Value = -2413.4546
$ws.Range("H112").Value = 4861
$ws.Range("J112").Value = 4861
$ws.Range("L112").Value = 14583
$ws.Range("N112").Value = -16799
$ws.Range("H122").Value = 1054.6154
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H128").Value = 94990
$ws.Range("J128").Value = 94990
$ws.Range("L128").Value = 94990
$ws.Range("N128").Value = -104950
$ws.Range("H132").Value = 1326.4509
$ws.Range("I132").Value = 784.375
$ws.Range("K132").Value = 2353.125
$ws.Range("M132").Value = 176.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 194.21428
$ws.Range("I5").Value = 102.375
$ws.Range("K5").Value = 102.375
$ws.Range("M5").Value = 9.625
$ws.Range("H32").Value = 38662.656
$ws.Range("I32").Value = 23306.375
$ws.Range("J32").Value = 112372.8
$ws.Range("K32").Value = 23306.375
$ws.Range("L32").Value = 112372.8
$ws.Range("M32").Value = -23019.375
$ws.Range("N32").Value = -112946.8
$ws.Range("H74").Value = 1657.1538
$ws.Range("I74").Value = 1404.174
$ws.Range("K74").Value = 1404.174
$ws.Range("M74").Value = -530.174
$ws.Range("H77").Value = 1657.1538
$ws.Range("I77").Value = 1404.174
$ws.Range("K77").Value = 7020.87
$ws.Range("M77").Value = -2652.87
$ws.Range("H97").Value = 3501.25
$ws.Range("I97").Value = 2332.6667
$ws.Range("K97").Value = 2332.6667
$ws.Range("M97").Value = -1836.6667
$ws.Range("H122").Value = 2778
$ws.Range("I122").Value = 2706.2856
$ws.Range("K122").Value = 8118.8568
$ws.Range("M122").Value = -5668.8568
$ws.Range("H132").Value = 2827.9285
$ws.Range("I132").Value = 1978.2941
$ws.Range("J132").Value = 4141
$ws.Range("K132").Value = 5934.8823
$ws.Range("L132").Value = 12423
$ws.Range("M132").Value = -3404.8823
$ws.Range("N132").Value = -17483

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 194.21428
$ws.Range("I4").Value = 102.375
$ws.Range("K4").Value = 102.375
$ws.Range("M4").Value = 12.625
$ws.Range("H134").Value = 1540.0476
$ws.Range("I134").Value = 1544.2632
$ws.Range("K134").Value = 4632.7896
$ws.Range("M134").Value = -2097.7896

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1705.4166
$ws.Range("I31").Value = 1477.0625
$ws.Range("J31").Value = 2162.125
$ws.Range("K31").Value = 1477.0625
$ws.Range("L31").Value = 2162.125
$ws.Range("M31").Value = -1182.0625
$ws.Range("N31").Value = -2752.125
$ws.Range("H34").Value = 1705.4166
$ws.Range("I34").Value = 1477.0625
$ws.Range("J34").Value = 2162.125
$ws.Range("K34").Value = 1477.0625
$ws.Range("L34").Value = 2162.125
$ws.Range("M34").Value = -1275.0625
$ws.Range("N34").Value = -2566.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8731097
$ws.Range("I4").Value = 3385412
$ws.Range("J4").Value = 45081750
$ws.Range("K4").Value = 10156236
$ws.Range("L4").Value = 135245250
$ws.Range("M4").Value = -10156124
$ws.Range("N4").Value = -135245474
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 5000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H18").Value = 830
$ws.Range("I18").Value = 830
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2490
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H44").Value = 36
$ws.Range("I44").Value = 36
$ws.Range("K44").Value = 108
$ws.Range("M44").Value = 290
$ws.Range("H46").Value = 168464.42
$ws.Range("I46").Value = 201172.3
$ws.Range("K46").Value = 603516.8999999999
$ws.Range("M46").Value = -603425.8999999999
$ws.Range("H75").Value = 8791.4
$ws.Range("J75").Value = 9029.538
$ws.Range("L75").Value = 27088.614
$ws.Range("N75").Value = -29084.614
$ws.Range("H78").Value = 8791.4
$ws.Range("J78").Value = 9029.538
$ws.Range("L78").Value = 81265.842
$ws.Range("N78").Value = -91249.842
$ws.Range("H129").Value = 120062.234
$ws.Range("J129").Value = 3055.3333
$ws.Range("L129").Value = 9165.999899999999
$ws.Range("N129").Value = -19165.9999
$ws.Range("H132").Value = 2424.12
$ws.Range("J132").Value = 2937.5
$ws.Range("L132").Value = 26437.5
$ws.Range("N132").Value = -31497.5
$ws.Range("H137").Value = 7696918
$ws.Range("J137").Value = 7168.857
$ws.Range("L137").Value = 21506.571
$ws.Range("N137").Value = -31706.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5736
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5736
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6276
$ws.Range("H73").Value = 5736
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5736
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7608
$ws.Range("H80").Value = 4374.5
$ws.Range("I80").Value = 3001.6667
$ws.Range("J80").Value = 5198.2
$ws.Range("K80").Value = 3001.6667
$ws.Range("L80").Value = 5198.2
$ws.Range("M80").Value = -2003.6667
$ws.Range("N80").Value = -7194.2
$ws.Range("H83").Value = 4374.5
$ws.Range("I83").Value = 3001.6667
$ws.Range("J83").Value = 5198.2
$ws.Range("K83").Value = 15008.3335
$ws.Range("L83").Value = 25991
$ws.Range("M83").Value = -10016.3335
$ws.Range("N83").Value = -35975
$ws.Range("H122").Value = 2699.4783
$ws.Range("I122").Value = 2915.5
$ws.Range("J122").Value = 2363.4443
$ws.Range("K122").Value = 8746.5
$ws.Range("L122").Value = 7090.3329
$ws.Range("M122").Value = -6296.5
$ws.Range("N122").Value = -11990.3329

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3463.9285
$ws.Range("I40").Value = 2541.6667
$ws.Range("K40").Value = 2541.6667
$ws.Range("M40").Value = -2405.6667
$ws.Range("H46").Value = 30294.4
$ws.Range("J46").Value = 2366.6667
$ws.Range("L46").Value = 2366.6667
$ws.Range("N46").Value = -2742.6667
$ws.Range("H76").Value = 10516.25
$ws.Range("J76").Value = 10516.25
$ws.Range("L76").Value = 10516.25
$ws.Range("N76").Value = -11192.25
$ws.Range("H79").Value = 10516.25
$ws.Range("J79").Value = 10516.25
$ws.Range("L79").Value = 10516.25
$ws.Range("N79").Value = -12856.25
$ws.Range("H82").Value = 993.3333
$ws.Range("I82").Value = 980
$ws.Range("K82").Value = 980
$ws.Range("M82").Value = -619
$ws.Range("H85").Value = 993.3333
$ws.Range("I85").Value = 980
$ws.Range("K85").Value = 980
$ws.Range("M85").Value = 268
$ws.Range("H100").Value = 27075.8
$ws.Range("I100").Value = 5197.25
$ws.Range("K100").Value = 5197.25
$ws.Range("M100").Value = -4656.25
$ws.Range("H136").Value = 3448.9524
$ws.Range("I136").Value = 3090.8235
$ws.Range("K136").Value = 9272.470499999999
$ws.Range("M136").Value = -6722.470499999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 44500
$ws.Range("J45").Value = 44000
$ws.Range("L45").Value = 44000
$ws.Range("N45").Value = -44982
$ws.Range("H68").Value = 25271
$ws.Range("J68").Value = 25271
$ws.Range("L68").Value = 25271
$ws.Range("N68").Value = -26893
$ws.Range("H71").Value = 25271
$ws.Range("J71").Value = 25271
$ws.Range("L71").Value = 75813
$ws.Range("N71").Value = -83925
$ws.Range("H122").Value = 718
$ws.Range("I122").Value = 664.6667
$ws.Range("K122").Value = 1994.0001
$ws.Range("M122").Value = 455.9999
$ws.Range("H132").Value = 8771
$ws.Range("I132").Value = 5963.0835
$ws.Range("K132").Value = 17889.2505
$ws.Range("M132").Value = -15359.2505
